# TC_5-FDR_E2E-FDR-2938-SYSTEM ADMIN REFERENCE DATA
#
# Rows 2 and 3 of the FDR_End_End_Receipt sheet had their "DrawDate" (col A)
# stored as a shared-string literal ("09/30/2020", t="s") instead of a real
# date serial. Re-entering the DrawDate as an actual date puts a numeric
# value in A2/A3 (which now carries the date through its existing date
# number format), which in turn recalculates the dependent PayDate formula
# in column B (=A+2). This also drops the now-unused "09/30/2020" entry
# from the shared-string table, renumbering later entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New DrawDate (10/12/2020, serial 44116) for both data rows; PayDate (col B,
# "=A+2" formula) recalculates automatically to 10/14/2020 (serial 44118).
$ws.Range("A2").Value = 44116
$ws.Range("A3").Value = 44116

# Leave the selection on B2, matching where the edit was made.
$ws.Range("B2").Select()
